# #5: property boat&car done
# Fix up the "汽車" (car) sheet header row and add the missing columns
# (property_category, category, date, legislator_name, legislator_id,
# source_file, index) that every other property sheet already carries.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("汽車")

# --- Row 1: proper column headers (previously this row mistakenly held
# a copy of the row-2 data values) ---
$ws.Range("B1").Value = "name"
$ws.Range("C1").Value = "capacity"
$ws.Range("D1").Value = "owner"
$ws.Range("E1").Value = "register_date"
$ws.Range("F1").Value = "register_reason"
$ws.Range("G1").Value = "acquire_value"
$ws.Range("H1").Value = "property_category"
$ws.Range("I1").Value = "category"
$ws.Range("J1").Value = "date"
$ws.Range("K1").Value = "legislator_name"
$ws.Range("L1").Value = "legislator_id"
$ws.Range("M1").Value = "source_file"
$ws.Range("N1").Value = "index"

# Match the bold/bordered header formatting already used by B1:G1
$ws.Range("B1:G1").Copy()
$ws.Range("H1:N1").PasteSpecial(-4122)

# --- Row 2: fill in the newly added data columns to match the pattern
# used on every other sheet (property_category/category/date/
# legislator_name/legislator_id/source_file/index) ---
$ws.Range("H2").Value = "land"
$ws.Range("I2").Value = "normal"
$ws.Range("J2").Value = "2013-12-11"
$ws.Range("K2").Value = "吳育仁"
$ws.Range("L2").Value = 1734
$ws.Range("M2").Value = "tmpbcc11"
$ws.Range("N2").Value = 39

# Match the formatting already used by B2:G2
$ws.Range("B2:G2").Copy()
$ws.Range("H2:N2").PasteSpecial(-4122)
